# Update "想去人数" (want-to-go count) figures in F column across sheets.
# Changes correspond to a refreshed data pull (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 826
$ws.Range("F3").Value  = 13639
$ws.Range("F4").Value  = 13441
$ws.Range("F5").Value  = 1043
$ws.Range("F6").Value  = 797
$ws.Range("F8").Value  = 585
$ws.Range("F10").Value = 16
$ws.Range("F11").Value = 38
$ws.Range("F12").Value = 734
$ws.Range("F13").Value = 2124
$ws.Range("F14").Value = 63
$ws.Range("F15").Value = 79
$ws.Range("F16").Value = 63
$ws.Range("F17").Value = 102
$ws.Range("F20").Value = 356
$ws.Range("F22").Value = 491
$ws.Range("F23").Value = 813

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 139
$ws.Range("F6").Value  = 152
$ws.Range("F7").Value  = 1332
$ws.Range("F10").Value = 53

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 89

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 826
$ws.Range("F4").Value  = 13639
$ws.Range("F5").Value  = 13441
$ws.Range("F6").Value  = 1043
$ws.Range("F7").Value  = 797
$ws.Range("F9").Value  = 585
$ws.Range("F11").Value = 16
$ws.Range("F12").Value = 38
$ws.Range("F13").Value = 734
$ws.Range("F16").Value = 2124
$ws.Range("F17").Value = 63
$ws.Range("F18").Value = 79
$ws.Range("F19").Value = 63
$ws.Range("F20").Value = 102
$ws.Range("F21").Value = 139
$ws.Range("F24").Value = 89
$ws.Range("F26").Value = 356
$ws.Range("F28").Value = 491
$ws.Range("F29").Value = 813
$ws.Range("F30").Value = 152
$ws.Range("F31").Value = 1332
$ws.Range("F35").Value = 53
